$d = $word.ActiveDocument

# Every target paragraph currently reads (across 2-3 runs, with mixed
# rPr/formatting):
#   ["<red-space> "]["Waktu Kampanye "]["2018 untuk Perseus: 30 Oktober-8
#   November dan 29 November-8 Desember"]
# (one of them is additionally prefixed, in the *same* <w:p>, by a
# "www.globeatnight.org" run + manual line break).
# Each becomes a single plain run (no rPr at all) reading:
#   "Waktu Kampanye Pegasus: 8-17 Oktober, 7-16 November,"
# Paragraph-level formatting (pStyle/jc/rPr) is untouched.

$oldMarker = "2018 untuk Perseus"
$newText   = "Waktu Kampanye Pegasus: 8-17 Oktober, 7-16 November,"

$newXmlTemplate = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr>' +
  '<w:pStyle w:val="BasicParagraph"/>' +
  '<w:jc w:val="center"/>' +
  '<w:rPr>' +
  '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Optima-Regular"/>' +
  '<w:color w:val="auto"/>' +
  '<w:sz w:val="26"/>' +
  '<w:szCs w:val="26"/>' +
  '</w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:t>{0}</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'
$newXml = $newXmlTemplate -f $newText

# Discover every paragraph that still carries the old text (rather than
# hard-coding indices) and fix them up back-to-front so rebuilding one
# paragraph can never disturb the index of one we haven't visited yet.
$targets = @()
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*$oldMarker*") {
        $targets += $i
    }
}

[array]::Reverse($targets)

foreach ($idx in $targets) {
    $p = $d.Paragraphs.Item($idx)
    $rng = $p.Range
    # Exclude the trailing paragraph mark so only the run content is wiped.
    $null = $rng.MoveEnd(1, -1)
    $null = $rng.Delete()
    $null = $rng.InsertXML($newXml)
}
